# 🔄 MAJ automatique BRVM via GitHub Actions
# Rename the existing sheet, then build a new "Top_YTD" sheet summarising the
# top 10 performers (by "Variation Totale (%)") already present on it.

$wb = $excel.ActiveWorkbook

# --- 1) Rename the original sheet ---------------------------------------
$wsMain = $wb.Worksheets.Item(1)
$wsMain.Name = "Recommandations"

# --- 2) Insert a new sheet right after it --------------------------------
$wsTop = $wb.Worksheets.Add($null, $wsMain)
$wsTop.Name = "Top_YTD"

# --- 3) Header row, reusing the same header style as the source sheet ---
$wsMain.Range("A1").Copy($wsTop.Range("A1"))
$wsMain.Range("D1").Copy($wsTop.Range("B1"))
$wsTop.Range("A1").Value = "Titre"
$wsTop.Range("B1").Value = "Progression YTD (%)"

# --- 4) Copy the Top 10 rows (already sorted descending on column D) ----
$topCount = 10
for ($i = 1; $i -le $topCount; $i++) {
    $srcRow = $i + 1
    $dstRow = $i + 1

    $title = $wsMain.Cells.Item($srcRow, 1).Value2
    $progression = $wsMain.Cells.Item($srcRow, 4).Value2

    $wsTop.Cells.Item($dstRow, 1).Value = $title
    $wsTop.Cells.Item($dstRow, 2).Value = $progression
}

# --- 5) Keep the original sheet as the active tab, like before the edit --
$wsMain.Select()
